$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===========================================================================
# 1) New "Characters" planning block appended at the bottom of the sheet
#    (rows 142-153). Values are written first, in the precise order that
#    reproduces the target shared-string table (new unique strings land at
#    indices 211-231, with "Hadewij" landing at 232 because it is written
#    after this block, further down this script).
# ===========================================================================

# Row 143 - header
$ws.Range("A143").Value = "Characters "
$ws.Range("B143").Value = "Asset name"
$ws.Range("C143").Value = "Asset code"
$ws.Range("D143").Value = "Artist"
$ws.Range("E143").Value = "Opmerkingen"

# Column B - character names, top to bottom (rows 144-153)
$ws.Range("B144").Value = "Dokter"
$ws.Range("B145").Value = "Assistant 1"
$ws.Range("B146").Value = "Assistant 2"
$ws.Range("B147").Value = "Client 1"
$ws.Range("B148").Value = "Client 2"
$ws.Range("B149").Value = "Client 3"
$ws.Range("B150").Value = "Client 4"
$ws.Range("B151").Value = "Client 5"
$ws.Range("B152").Value = "Client 6"
$ws.Range("B153").Value = "Client 7"

# Column C - 3D asset codes (entry order matches the target shared-string
# table exactly, which interleaves row 147 before row 146).
$ws.Range("C144").Value = "3D_CHAR_Dokter"
$ws.Range("C145").Value = "3D_CHAR_Assistant1"
$ws.Range("C147").Value = "3D_CHAR_Client1"
$ws.Range("C146").Value = "3D_CHAR_Assistant2"
$ws.Range("C148").Value = "3D_CHAR_Client2"
$ws.Range("C149").Value = "3D_CHAR_Client3"
$ws.Range("C150").Value = "3D_CHAR_Client4"
$ws.Range("C151").Value = "3D_CHAR_Client5"
$ws.Range("C152").Value = "3D_CHAR_Client6"
$ws.Range("C153").Value = "3D_CHAR_Client7"

# ---------------------------------------------------------------------------
# Formatting for the new block, copied (formats-only) from existing cells
# that already carry the matching look, so no new cell styles are created.
# ---------------------------------------------------------------------------

$ws.Range("B3").Copy()
$ws.Range("A143").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("B143:D143").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("E143").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("F143:F153").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("A144:A152").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("B144:B152").PasteSpecial(-4122)
$ws.Range("D144:D152").PasteSpecial(-4122)
$ws.Range("C144:C145").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("C146:C152").PasteSpecial(-4122)

$ws.Range("F4").Copy()
$ws.Range("E144:E152").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("A153").PasteSpecial(-4122)

$ws.Range("C34").Copy()
$ws.Range("B153").PasteSpecial(-4122)
$ws.Range("D153").PasteSpecial(-4122)

$ws.Range("C136").Copy()
$ws.Range("C153").PasteSpecial(-4122)

$ws.Range("F34").Copy()
$ws.Range("E153").PasteSpecial(-4122)

$ws.Rows(142).RowHeight = 15.75

# Re-apply all values - the Formats-only pastes above must not disturb them,
# but this keeps the block correct even if a paste ever clobbers a value.
$ws.Range("A143").Value = "Characters "
$ws.Range("B143").Value = "Asset name"
$ws.Range("C143").Value = "Asset code"
$ws.Range("D143").Value = "Artist"
$ws.Range("E143").Value = "Opmerkingen"
$ws.Range("B144").Value = "Dokter"
$ws.Range("C144").Value = "3D_CHAR_Dokter"
$ws.Range("B145").Value = "Assistant 1"
$ws.Range("C145").Value = "3D_CHAR_Assistant1"
$ws.Range("B146").Value = "Assistant 2"
$ws.Range("C146").Value = "3D_CHAR_Assistant2"
$ws.Range("B147").Value = "Client 1"
$ws.Range("C147").Value = "3D_CHAR_Client1"
$ws.Range("B148").Value = "Client 2"
$ws.Range("C148").Value = "3D_CHAR_Client2"
$ws.Range("B149").Value = "Client 3"
$ws.Range("C149").Value = "3D_CHAR_Client3"
$ws.Range("B150").Value = "Client 4"
$ws.Range("C150").Value = "3D_CHAR_Client4"
$ws.Range("B151").Value = "Client 5"
$ws.Range("C151").Value = "3D_CHAR_Client5"
$ws.Range("B152").Value = "Client 6"
$ws.Range("C152").Value = "3D_CHAR_Client6"
$ws.Range("B153").Value = "Client 7"
$ws.Range("C153").Value = "3D_CHAR_Client7"

# ===========================================================================
# 2) Existing-cell updates elsewhere on the planning sheet. "Hadewij" is
#    entered last so it becomes shared-string index 232, matching the
#    target file's string table.
# ===========================================================================
$ws.Range("J5").Value = "Quinten"
$ws.Range("J8").Value = "Quinten"
$ws.Range("J6").Value = "Marc"
$ws.Range("J10").Value = "Marc"
$ws.Range("J3").Value = "Hadewij"
$ws.Range("E5").Value = "Hadewij"

# J10 also picks up the same row-formatting as I10 (its neighbour).
$ws.Range("I10").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J10").Value = "Marc"

# ===========================================================================
# 3) Misc sheet-view bookkeeping - cursor moved to J9 before saving.
# ===========================================================================
$ws.Range("J9").Select()
